$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddVacancy")

# Update the sample "add vacancy" test data row: vacancyName / jobTitle / hiringManager
$ws.Range("C2").Value = "Test3"
$ws.Range("D2").Value = "Software developer"
$ws.Range("E2").Value = "Tuong  Ma"

# Widen the newly meaningful columns to fit their content (mirrors the other sheets'
# auto-fitted columns) and select/activate this sheet as the new active tab.
$ws.Columns.Item(1).ColumnWidth = 8.166666666666666
$ws.Columns.Item(4).ColumnWidth = 12.608072916666666
$ws.Columns.Item(5).ColumnWidth = 11.830729166666666

$ws.Range("E2").Select()
